# Logboek-Noël.xlsx
# Commit: "C# en Noël logboek" - append two more logbook days (Week 3:
# Thursday 11-05-2017 and Friday 12-05-2017, "C# app gewerkt") and move the
# viewport/selection further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Thursday 11-05-2017 (serial 42866) -> rows 46/47, mirrors rows 43/44 ---
$ws.Range("A43").Copy()
$ws.Range("A46").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A46").Value = 42866

$ws.Range("E43").Copy()
$ws.Range("E46").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("B46").Value = "Week 3"
$ws.Range("C46").Value = "do"
$ws.Range("G46").Value = "School"
$ws.Range("I46").Value = "C# app gewerkt"

$ws.Range("B47").Value = "Week 3"
$ws.Range("G47").Value = "School"

# --- Friday 12-05-2017 (serial 42867) -> rows 49/50, mirrors rows 43/44 ---
$ws.Range("A43").Copy()
$ws.Range("A49").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A49").Value = 42867

$ws.Range("E43").Copy()
$ws.Range("E49").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("B49").Value = "Week 3"
$ws.Range("C49").Value = "vr"
$ws.Range("G49").Value = "School"
$ws.Range("I49").Value = "C# app gewerkt"

$ws.Range("B50").Value = "Week 3"
$ws.Range("G50").Value = "School"

$excel.CutCopyMode = 0

# --- Scroll the view down and move the selection, like the author did ---
$null = $ws.Range("E53").Select()
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
